# Add three new rows to the "Dic_Disagg_Ausprägungen" dictionary sheet so
# that the alphabetically-sorted list of codes (column A) gains:
#   A_AIRPOLL_TOTAL      (after A_AIRPOLL_SO2,   before A_AREA_EU)
#   A_SERIES_AUDIOVIS    (after A_SERIES_ANNUALVAL, before A_SERIES_BEH)
#   A_SERIES_KULTUR      (after A_SERIES_IMP,    before A_SERIES_LOCATIONS)
# Every existing row beneath each insertion point shifts down by one,
# growing the sheet from 167 to 170 data/header rows (dimension A1:D170).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Insert "A_AIRPOLL_TOTAL" as the new row 27.
# ---------------------------------------------------------------------
$ws.Rows("27:27").Insert()
$ws.Range("A27").Value = "A_AIRPOLL_TOTAL"
$ws.Range("B27").Value = "K_AIRPOLL"
$ws.Range("C27").Value = "Insgesamt"
$ws.Range("D27").Value = "XXXInsgesamt"
# Copy formatting from the row directly above so the new row matches the
# sheet's standard data-row style instead of Excel's inserted default.
$ws.Range("A26:D26").Copy()
$ws.Range("A27:D27").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 2) Insert "A_SERIES_AUDIOVIS" as the new row 89 (old row 88 + 1 shift).
# ---------------------------------------------------------------------
$ws.Rows("89:89").Insert()
$ws.Range("A89").Value = "A_SERIES_AUDIOVIS"
$ws.Range("B89").Value = "K_SERIES"
$ws.Range("C89").Value = "Produktionen audiovisueller Medien, bei denen ökologische Standards eingehalten wurden"
$ws.Range("D89").Value = "XXXProduktionen audiovisueller Medien, bei denen ökologische Standards eingehalten wurden"
$ws.Range("A88:D88").Copy()
$ws.Range("A89:D89").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 3) Insert "A_SERIES_KULTUR" as the new row 118 (old row 117 + 2 shift).
# ---------------------------------------------------------------------
$ws.Rows("118:118").Insert()
$ws.Range("A118").Value = "A_SERIES_KULTUR"
$ws.Range("B118").Value = "K_SERIES"
$ws.Range("C118").Value = "Nachhaltigkeitszertifizierte Kultur- und Medieneinrichtungen"
$ws.Range("D118").Value = "XXXNachhaltigkeitszertifizierte Kultur- und Medieneinrichtungen"
$ws.Range("A117:D117").Copy()
$ws.Range("A118:D118").PasteSpecial(-4122)
$excel.CutCopyMode = 0
